# Commit: "Update set from data objects"
#
# sheet1 keeps its header row, but the 3 data rows shift down by one row
# (a blank spacer row is inserted right under the header) and the
# trailing blank rows/columns that only existed because the sheet's
# used-range was originally padded out to J10 are trimmed away. A new
# sheet2 is added next to it holding a trimmed copy of the original
# layout (header + 3 data rows, column E emptied out).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- sheet1 -----------------------------------------------------------

# shift the 3 data rows down by one row (values only)
$ws1.Range("A2:E4").Copy()
$ws1.Range("A3").PasteSpecial(-4163)
$ws1.Application.CutCopyMode = $false

# carry the highlighted-cell formatting from C2 to its new home at C3
$ws1.Range("C2").Copy()
$ws1.Range("C3").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

# row 2 becomes an empty spacer row: clear its old values/format, then
# stamp it with the plain header formatting so the cells stay present
# (as empty placeholders) instead of disappearing entirely
$ws1.Range("A2:E2").ClearContents()
$ws1.Range("C2").ClearFormats()
$ws1.Range("A1:E1").Copy()
$ws1.Range("A2").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

# drop the unused trailing blank rows/columns
$ws1.Columns("F:J").Delete()
$ws1.Rows("6:10").Delete()

# cosmetics carried by the edit
$ws1.Range("C1").ColumnWidth = 99.15
$ws1.Application.ActiveWindow.DisplayGridlines = $false

# --- sheet2: new sheet with a trimmed copy of the original data ------

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "sheet2"

$ws2.Range("A1").Value = "A"
$ws2.Range("B1").Value = "B"
$ws2.Range("C1").Value = "C"
$ws2.Range("D1").Value = "D"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = 3
$ws2.Range("D2").Value = 4

$ws2.Range("A3").Value = 6
$ws2.Range("B3").Value = 7
$ws2.Range("C3").Value = 8
$ws2.Range("D3").Value = 9

$ws2.Range("A4").Value = 11
$ws2.Range("B4").Value = 12
$ws2.Range("C4").Value = 13
$ws2.Range("D4").Value = 14

# give C2 the same highlighted-cell formatting as sheet1's C3, and make
# E1:E4 present as empty placeholders (formats-only paste from the
# plain header cells) to match the original column count
$ws1.Range("C3").Copy()
$ws2.Range("C2").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

$ws2.Range("A1:A4").Copy()
$ws2.Range("E1:E4").PasteSpecial(-4122)
$ws2.Application.CutCopyMode = $false

$ws2.Range("C1").ColumnWidth = 99.15
$ws2.Application.ActiveWindow.DisplayGridlines = $false
